# Applies the "Excel with test data modified" commit:
#  - Info sheet: move the selection from A1:A18 to B6 (no other changes)
#  - Liste_ChatGPT sheet: scroll the frozen pane back up to A2 (selection stays A41:Q41)
#  - Mapping_JSON sheet: becomes a 3-column mapping table
#       column B is relabelled to hold the lower-case JSON attribute names
#       a new column C is added holding the upper-case JSON attribute names
#       header B1 becomes "Attributname JSON (klein)", new header C1 "Attributname JSON (groß)"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Info sheet: selection changes from A1:A18 to B6
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("B6").Select()

# ---------------------------------------------------------------------------
# 2) Mapping_JSON sheet: rework into a 3 column table
# ---------------------------------------------------------------------------
$wsMap = $wb.Worksheets.Item("Mapping_JSON")
$wsMap.Activate()

# New lower-case JSON attribute names for column B (rows 1-18)
$colB = @(
  "Attributname JSON (klein)",
  "sst_name",
  "sst_version",
  "sst_valid_from",
  "sst_valid_to",
  "sst_is_active",
  "sst_responsible",
  "sst_auth_method",
  "sst_authorization",
  "sst_crypto",
  "sst_format",
  "sst_trigger",
  "sst_protocol",
  "sst_transport",
  "partner_1_system",
  "partner_1_responsible",
  "partner_2_system",
  "partner_2_responsible"
)

# New upper-case JSON attribute names for the new column C (rows 1-18)
$colC = @(
  "Attributname JSON (groß)",
  "SST_NAME",
  "SST_VERSION",
  "SST_VALID_FROM",
  "SST_VALID_TO",
  "SST_IS_ACTIVE",
  "SST_RESPONSIBLE",
  "SST_AUTH_METHOD",
  "SST_AUTHORIZATION",
  "SST_CRYPTO",
  "SST_FORMAT",
  "SST_TRIGGER",
  "SST_PROTOCOL",
  "SST_TRANSPORT",
  "PARTNER_1_SYSTEM",
  "PARTNER_1_RESPONSIBLE",
  "PARTNER_2_SYSTEM",
  "PARTNER_2_RESPONSIBLE"
)

for ($i = 0; $i -lt $colB.Length; $i++) {
  $row = $i + 1
  $wsMap.Cells.Item($row, 2).Value = $colB[$i]
  $wsMap.Cells.Item($row, 3).Value = $colC[$i]
}

# Give column C the same look (fill/border/alignment) as column B
$wsMap.Range("B1").Copy()
$wsMap.Range("C1").PasteSpecial(-4122)
$wsMap.Range("B2:B18").Copy()
$wsMap.Range("C2:C18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the new column's width to the existing ones
$wsMap.Columns.Item(3).ColumnWidth = 29.8

# Page setup now mirrors the other sheets (portrait, paper size 9)
$wsMap.PageSetup.PaperSize = 9
$wsMap.PageSetup.Orientation = 1

# Selection on Mapping_JSON ends up on A3 after the edits
$wsMap.Range("A3").Select()

# ---------------------------------------------------------------------------
# 3) Liste_ChatGPT sheet: scroll frozen pane back to the top (A2)
# ---------------------------------------------------------------------------
$wsList = $wb.Worksheets.Item("Liste_ChatGPT")
$wsList.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$wsList.Range("A41:Q41").Select()

# Leave Mapping_JSON as the active/visible sheet (matches activeTab=1 / tabSelected)
$wsMap.Activate()
